$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.033231139183044
$ws.Range("B1").Value = 1.447754502296448
$ws.Range("C1").Value = 3.806950330734253
$ws.Range("D1").Value = 2.118509769439697
$ws.Range("E1").Value = 0.8314157128334045
